# The underlying data rows (2-17) got re-sorted/re-ordered by the source
# system: row 3's record moved to the top (new row 2), row 17's record moved
# right under it (new row 3), the old row 2 record dropped to new row 4, and
# every other record (old rows 4-16) simply shifted down by one row (into
# new rows 5-17). Only a handful of columns actually carry per-record data
# (A = Id, M/N = Aktivitet/Metod, Q/R = Ost/Nord coordinates, AC = Publik
# kommentar) - everything else on the row is identical for every record, so
# we only need to move those six columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17

# New-row -> old-row mapping that reproduces the observed re-order.
$mapping = @{}
$mapping[2] = 3
$mapping[3] = 17
$mapping[4] = 2
for ($r = 5; $r -le $lastRow; $r++) { $mapping[$r] = $r - 1 }

$cols = @("A", "M", "N", "Q", "R", "AC")

# 1) Snapshot the current ("old") values for every tracked column, per row,
#    before anything is overwritten.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rng = $ws.Range($col + $r)
        $rowVals[$col] = @{ Value = $rng.Value2; Empty = ($rng.Text -eq "") }
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the snapshotted values back out in the new order. Cells that were
#    genuinely blank in the source row are cleared (rather than written with
#    an empty string) so a column that didn't apply to that record stays
#    blank in the destination too.
foreach ($newRow in ($mapping.Keys | Sort-Object)) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($col in $cols) {
        $dest = $ws.Range($col + $newRow)
        $info = $src[$col]
        if ($info.Empty) {
            $dest.ClearContents()
        } else {
            $dest.Value2 = $info.Value
        }
    }
}
